$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 254 (pushing old rows 254-345 down to 256-347)
$ws.Rows("254:255").Insert()

# New row 254 data
$ws.Range("A254").Value = 10
$ws.Range("B254").Value = "Vega Modelo de Temuco"
$ws.Range("C254").Value = "La Araucanía"
$ws.Range("D254").Value = 44559
$ws.Range("E254").Value = 9
$ws.Range("F254").Value = 100112023
$ws.Range("G254").Value = "Brócoli"
$ws.Range("H254").Value = "Sin especificar"
$ws.Range("I254").Value = "Primera"
$ws.Range("J254").Value = 650
$ws.Range("K254").Value = 800
$ws.Range("L254").Value = 800
$ws.Range("M254").Value = 800
$ws.Range("N254").Value = "$/unidad"
$ws.Range("O254").Value = "Región de O'Higgins"
$ws.Range("P254").Value = 800
$ws.Range("Q254").Value = 1
$ws.Range("R254").Value = "Hortaliza"

# New row 255 data
$ws.Range("A255").Value = 10
$ws.Range("B255").Value = "Vega Modelo de Temuco"
$ws.Range("C255").Value = "La Araucanía"
$ws.Range("D255").Value = 44559
$ws.Range("E255").Value = 9
$ws.Range("F255").Value = 100112023
$ws.Range("G255").Value = "Brócoli"
$ws.Range("H255").Value = "Sin especificar"
$ws.Range("I255").Value = "Primera"
$ws.Range("J255").Value = 1250
$ws.Range("K255").Value = 800
$ws.Range("L255").Value = 800
$ws.Range("M255").Value = 800
$ws.Range("N255").Value = "$/unidad"
$ws.Range("O255").Value = "Región del Maule"
$ws.Range("P255").Value = 800
$ws.Range("Q255").Value = 1
$ws.Range("R255").Value = "Hortaliza"
